# Applies re-ordering of the semantic matching results rows (A2:C28).
# The underlying (Unmatched Item, Matched Master Item, Similarity Score)
# triples are unchanged as a set; only their row order changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(" add fresh paint", "Desk for Repaint", 0.5055999755859375),
    @(" air con grill is loose", "Air con Grill Loose", 0.9781000018119812),
    @("Heineken", "Sweetener Sachet", 0.3255999982357025),
    @("  gas burner", "Air conditioner is faulty", 0.4805000126361847),
    @("  hair dryer", "Hair dryer not working", 0.9121000170707703),
    @("Veuve Clicquot Brut NV  ", "Loofah", 0.2937999963760376),
    @("  fan body needs paint job", "Ceiling Painting", 0.4948999881744385),
    @(" bed ceiling repaint", "Bedroom Ceiling Light Flickering", 0.6807000041007996),
    @("  fan vent in kitchen", "Exhaust Fan", 0.6567999720573425),
    @(" AC duct border to repaint", "Desk for Repaint", 0.5903000235557556),
    @("  no  dnd light showing", "Light Dimmer Not Functioning", 0.5509999990463257),
    @("BT", "Express laundry service", 0.2673999965190887),
    @(" add grout to tile corner", "Entrance Wall Socket", 0.4327999949455261),
    @("peppermint tea", "Chamomile tea", 0.6674000024795532),
    @("Writing Table Light Fused", "Writing Table Light Fused", 1),
    @(" blu tooth liight", "Dental kit", 0.5746999979019165),
    @(" door glass strip is broken", "Glass Wall - Broken or Cracked", 0.7461000084877014),
    @("set up the bait rail", "Tape to seal box", 0.4758999943733215),
    @("  pond vacuum", "Pond Water Low Level", 0.6601999998092651),
    @("  my basin tap is loose", "Basin Tap Loose", 0.9185000061988831),
    @("  no space bar counter", "Bar Counter High Chair Defective", 0.6266000270843506),
    @("bathroom crumb polish", "Cleaning of Bathroom", 0.6427000164985657),
    @("  my bag rack is loose", "Trash Bag", 0.5766000151634216),
    @("  darjeeling tebags", "Tea bags", 0.4634000062942505),
    @(" air band dirty", "Air con Grill Dirty", 0.6209999918937683),
    @(" buff gel", "Bath gel", 0.6707999706268311),
    @("  main lobby cleaning", "2nd Floor Guest Lift Lobby to Clean", 0.6998999714851379)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
